$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E are treated as text so numeric-looking values
# (e.g. "64.308.43", "6.07") are not coerced into numbers/dates by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '64.308.43'
$ws.Range("E2").Value = '  -2.50%  '

# Row 3
$ws.Range("D3").Value = '3.177.20'
$ws.Range("E3").Value = '  -7.76%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").Value = '563.46'
$ws.Range("E5").Value = '  -3.67%  '

# Row 6
$ws.Range("D6").Value = '170.88'
$ws.Range("E6").Value = '  -1.76%  '

# Row 7
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value = '0.610'
$ws.Range("E7").Value = '  +1.21%  '

# Row 8
$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.00%  '

# Row 9
$ws.Range("D9").Value = '3.171.68'
$ws.Range("E9").Value = '  -7.88%  '

# Row 10
$ws.Range("E10").Value = '  -5.64%  '

# Row 11
$ws.Range("D11").Value = '6.62'
$ws.Range("E11").Value = '  -4.94%  '

# Row 12
$ws.Range("D12").Value = '0.396'
$ws.Range("E12").Value = '  -4.44%  '

# Row 13
$ws.Range("D13").Value = '3.722.25'
$ws.Range("E13").Value = '  -8.04%  '

# Row 14
$ws.Range("D14").Value = '0.136'
$ws.Range("E14").Value = '  +1.15%  '

# Row 15
$ws.Range("D15").Value = '27.42'
$ws.Range("E15").Value = '  -5.89%  '

# Row 16
$ws.Range("D16").Value = '64.309.60'
$ws.Range("E16").Value = '  -2.54%  '

# Row 17
$ws.Range("E17").Value = '  -5.03%  '

# Row 18
$ws.Range("D18").Value = '3.175.86'
$ws.Range("E18").Value = '  -7.73%  '

# Row 19
$ws.Range("D19").Value = '5.71'
$ws.Range("E19").Value = '  -4.12%  '

# Row 20
$ws.Range("D20").Value = '13.05'
$ws.Range("E20").Value = '  -5.62%  '

# Row 21
$ws.Range("D21").Value = '352.88'
$ws.Range("E21").Value = '  -4.60%  '

# Row 22
$ws.Range("E22").Value = '  -5.13%  '

# Row 23
$ws.Range("E23").Value = '  +0.03%  '

# Row 24
$ws.Range("D24").Value = '69.19'
$ws.Range("E24").Value = '  -4.29%  '

# Row 25
$ws.Range("E25").Value = '  -5.08%  '

# Row 26
$ws.Range("E26").Value = '  -2.90%  '

# Row 27
$ws.Range("E27").Value = '  -1.16%  '

# Row 28
$ws.Range("E28").Value = '  -1.76%  '

# Row 29
$ws.Range("E29").Value = '  +0.04%  '

# Row 30
$ws.Range("D30").Value = '5.64'
$ws.Range("E30").Value = '  -2.71%  '

# Row 31
$ws.Range("E31").Value = '  -0.16%  '

# Row 32
$ws.Range("E32").Value = '  -4.38%  '

# Row 33
$ws.Range("D33").Value = '22.14'
$ws.Range("E33").Value = '  -6.24%  '

# Row 34
$ws.Range("D34").Value = '6.67'
$ws.Range("E34").Value = '  -4.91%  '

# Row 35
$ws.Range("E35").Value = '  -5.62%  '

# Row 36
$ws.Range("E36").Value = '  -6.00%  '

# Row 37
$ws.Range("D37").Value = '154.61'
$ws.Range("E37").Value = '  -4.33%  '

# Row 38
$ws.Range("D38").Value = '0.810'
$ws.Range("E38").Value = '  -7.92%  '

# Row 39
$ws.Range("D39").Value = '26.04'
$ws.Range("E39").Value = '  -8.14%  '

# Row 40
$ws.Range("E40").Value = '  -2.27%  '

# Row 41
$ws.Range("E41").Value = '  -5.54%  '

# Row 42
$ws.Range("D42").Value = '2.631.82'
$ws.Range("E42").Value = '  -5.38%  '

# Row 43
$ws.Range("D43").Value = '4.18'
$ws.Range("E43").Value = '  -6.20%  '

# Row 44
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '6.07'
$ws.Range("E44").Value = '  -6.01%  '

# Row 45
$ws.Range("B45").Value = 'Bittensor'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D45").Value = '332.04'
$ws.Range("E45").Value = '  +1.11%  '

# Row 46
$ws.Range("E46").Value = '  -4.44%  '

# Row 47
$ws.Range("E47").Value = '  -2.55%  '

# Row 48
$ws.Range("D48").Value = '23.89'
$ws.Range("E48").Value = '  -4.24%  '

# Row 49
$ws.Range("D49").Value = '0.0271'
$ws.Range("E49").Value = '  -7.08%  '

# Row 50
$ws.Range("E50").Value = '  -1.10%  '

# Row 51
$ws.Range("D51").Value = '0.999'
$ws.Range("E51").Value = '  -0.03%  '

# Restore default (no explicit number format) style for D:E so the
# saved XML does not carry a spurious style index on these cells.
$ws.Range("D2:E51").Style = "Normal"
